# "Se agregaron nuevos tempos" - new time-tracking entries were logged
# into the "Casos de Uso" sheet's day-by-day consumption columns, the
# status of one task moved from "Por iniciar" to "En proceso", and the
# one-off bold/underline highlight on B14 was reverted back to the
# normal task-name style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# --- new hours consumed ("Cons.") entered for several tasks/days ---
# Row 13 ("GUI_CU01..." task): 7 hours consumed on Day 2 ("Cons." column K)
$ws.Range("K13").Value = 7

# Row 14 ("GUI_CU02" task): 23 hours consumed on Day 1 ("Cons." column H)
$ws.Range("H14").Value = 23

# Row 15 ("GUI_CU03" task): 11 hours consumed on Day 1 ("Cons." column H)
$ws.Range("H15").Value = 11

# Row 16: 44 hours consumed on Day 3 ("Cons." column N)
$ws.Range("N16").Value = 44

# Row 18: 48 hours consumed on Day 2 ("Cons." column K)
$ws.Range("K18").Value = 48

# --- row 14 status moves from "Por iniciar" to "En proceso" ---
# copy the formatting already used for "En proceso" status cells (e.g. F13)
# onto F14, then set the new status text
$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = "En proceso"

# --- B14 loses its one-off bold+underline highlight, back to normal ---
# copy the normal task-name formatting from a neighboring row (B13) onto B14
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)

# restore the selection so the active cell matches where the edit left off
$ws.Range("F16").Select()

$excel.CutCopyMode = $false
